# NYPD CompStat weekly report refresh: new crime data collected.
# Updates the "Volume/Number" and "Report Covering the Week" header text,
# plus the full crime-complaint statistics grid (rows 15-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text (rich-text shared strings) — edit from the back of each
# string forward so earlier character offsets aren't shifted by a change
# in string length later in the same cell.
# ---------------------------------------------------------------------

# A8: "Volume 30   Number  36" -> "...  37"
$ws.Range("A8").Characters(21, 2).Text = "37"

# C9: "Report Covering the Week  9/4/2023  Through  9/10/2023"
#   -> "...  9/11/2023  Through  9/17/2023"
$ws.Range("C9").Characters(46, 9).Text = "9/17/2023"
$ws.Range("C9").Characters(27, 8).Text = "9/11/2023"

# ---------------------------------------------------------------------
# Helper: PasteSpecial-formats-only paste code (xlPasteFormats) so that
# when a cell flips between a numeric style and the "N/A" text style, it
# adopts the destination style's existing xf record instead of minting a
# brand-new one.
# ---------------------------------------------------------------------
$xlPasteFormats = -4122

function Set-NumFromText($row, $col, $styleSourceRow, $val) {
    $ws.Cells.Item($styleSourceRow, $col).Copy() | Out-Null
    $ws.Cells.Item($row, $col).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Cells.Item($row, $col).Value = $val
}

function Set-TextFromNum($row, $col, $text) {
    $ws.Cells.Item($row, $col).NumberFormat = "@"
    $ws.Cells.Item($row, $col).Value = $text
}

# ---------------------------------------------------------------------
# Row 15 — Rape
# ---------------------------------------------------------------------
$ws.Cells.Item(15, 7).Value = 3            # G15
$ws.Cells.Item(15, 10).Value = 14          # J15
$ws.Cells.Item(15, 11).Value = -42.857142857142   # K15
$ws.Cells.Item(15, 12).Value = -57.894736842105   # L15

# ---------------------------------------------------------------------
# Row 16 — Robbery
# ---------------------------------------------------------------------
$ws.Cells.Item(16, 3).Value = 8            # C16
$ws.Cells.Item(16, 4).Value = 4            # D16
$ws.Cells.Item(16, 5).Value = 100          # E16
$ws.Cells.Item(16, 6).Value = 23           # F16
$ws.Cells.Item(16, 7).Value = 19           # G16
$ws.Cells.Item(16, 8).Value = 21.052631578947     # H16
$ws.Cells.Item(16, 9).Value = 125          # I16
$ws.Cells.Item(16, 10).Value = 137         # J16
$ws.Cells.Item(16, 11).Value = -8.759124087591    # K16
$ws.Cells.Item(16, 12).Value = -3.846153846153    # L16
$ws.Cells.Item(16, 13).Value = 0.806451612903     # M16
$ws.Cells.Item(16, 14).Value = -83.681462140992   # N16

# ---------------------------------------------------------------------
# Row 17 — Fel. Assault
# ---------------------------------------------------------------------
$ws.Cells.Item(17, 4).Value = 5            # D17
$ws.Cells.Item(17, 5).Value = -20          # E17
$ws.Cells.Item(17, 6).Value = 17           # F17
$ws.Cells.Item(17, 8).Value = 0            # H17
$ws.Cells.Item(17, 9).Value = 161          # I17
$ws.Cells.Item(17, 10).Value = 171         # J17
$ws.Cells.Item(17, 11).Value = -5.847953216374    # K17
$ws.Cells.Item(17, 12).Value = 8.053691275167     # L17
$ws.Cells.Item(17, 13).Value = 69.473684210526    # M17
$ws.Cells.Item(17, 14).Value = -20.689655172413   # N17

# ---------------------------------------------------------------------
# Row 18 — Burglary  (D18, E18 flip from numbers to the "N/A" text pair)
# ---------------------------------------------------------------------
Set-TextFromNum 18 4 "0"        # D18 -> "0"
Set-TextFromNum 18 5 "***.*"    # E18 -> "***.*"
$ws.Cells.Item(18, 6).Value = 17           # F18
$ws.Cells.Item(18, 7).Value = 21           # G18
$ws.Cells.Item(18, 8).Value = -19.047619047619    # H18
$ws.Cells.Item(18, 9).Value = 161          # I18
$ws.Cells.Item(18, 11).Value = -22.596153846153   # K18
$ws.Cells.Item(18, 12).Value = -0.617283950617    # L18
$ws.Cells.Item(18, 13).Value = -14.814814814814   # M18
$ws.Cells.Item(18, 14).Value = -85.839929639401   # N18

# ---------------------------------------------------------------------
# Row 19 — Gr. Larceny
# ---------------------------------------------------------------------
$ws.Cells.Item(19, 3).Value = 24           # C19
$ws.Cells.Item(19, 4).Value = 21           # D19
$ws.Cells.Item(19, 5).Value = 14.285714285714     # E19
$ws.Cells.Item(19, 6).Value = 79           # F19
$ws.Cells.Item(19, 7).Value = 82           # G19
$ws.Cells.Item(19, 8).Value = -3.658536585365     # H19
$ws.Cells.Item(19, 9).Value = 724          # I19
$ws.Cells.Item(19, 10).Value = 719         # J19
$ws.Cells.Item(19, 11).Value = 0.695410292072     # K19
$ws.Cells.Item(19, 12).Value = 39.768339768339    # L19
$ws.Cells.Item(19, 13).Value = -28.316831683168   # M19
$ws.Cells.Item(19, 14).Value = -62.291666666666   # N19

# ---------------------------------------------------------------------
# Row 20 — G.L.A.
# ---------------------------------------------------------------------
$ws.Cells.Item(20, 3).Value = 2            # C20
$ws.Cells.Item(20, 5).Value = 0            # E20
$ws.Cells.Item(20, 7).Value = 4            # G20
$ws.Cells.Item(20, 8).Value = 25           # H20
$ws.Cells.Item(20, 9).Value = 47           # I20
$ws.Cells.Item(20, 10).Value = 55          # J20
$ws.Cells.Item(20, 11).Value = -14.545454545454   # K20
$ws.Cells.Item(20, 12).Value = 11.904761904761    # L20
$ws.Cells.Item(20, 13).Value = 56.666666666666    # M20
$ws.Cells.Item(20, 14).Value = -94.515752625437   # N20

# ---------------------------------------------------------------------
# Row 21 — TOTAL
# ---------------------------------------------------------------------
$ws.Cells.Item(21, 3).Value = 41           # C21
$ws.Cells.Item(21, 4).Value = 33           # D21
$ws.Cells.Item(21, 5).Value = 24.242424242424    # E21
$ws.Cells.Item(21, 6).Value = 141          # F21
$ws.Cells.Item(21, 7).Value = 148          # G21
$ws.Cells.Item(21, 8).Value = -4.729729729729    # H21
$ws.Cells.Item(21, 9).Value = 1227         # I21
$ws.Cells.Item(21, 10).Value = 1306        # J21
$ws.Cells.Item(21, 11).Value = -6.049004594180   # K21
$ws.Cells.Item(21, 12).Value = 20.176297747306   # L21
$ws.Cells.Item(21, 13).Value = -15.728021978022  # M21
$ws.Cells.Item(21, 14).Value = -74.964293001428  # N21

# ---------------------------------------------------------------------
# Row 22 — Transit
# ---------------------------------------------------------------------
$ws.Cells.Item(22, 7).Value = 2            # G22
$ws.Cells.Item(22, 8).Value = 100          # H22
$ws.Cells.Item(22, 9).Value = 62           # I22
$ws.Cells.Item(22, 11).Value = -3.125      # K22
$ws.Cells.Item(22, 12).Value = 58.974358974359   # L22
$ws.Cells.Item(22, 13).Value = 16.981132075471   # M22

# ---------------------------------------------------------------------
# Row 23 — Housing  (D23/E23/G23/H23 flip from "N/A" text to numbers;
# borrow row 24's number styles so the xf record is reused, not minted)
# ---------------------------------------------------------------------
Set-NumFromText 23 4 24 1          # D23 = 1
Set-NumFromText 23 5 24 -100       # E23 = -100
Set-NumFromText 23 7 24 1          # G23 = 1
Set-NumFromText 23 8 24 -100       # H23 = -100
$ws.Cells.Item(23, 10).Value = 14          # J23
$ws.Cells.Item(23, 11).Value = -50         # K23
$ws.Cells.Item(23, 12).Value = -12.5       # L23

# ---------------------------------------------------------------------
# Row 24 — Petit Larceny
# ---------------------------------------------------------------------
$ws.Cells.Item(24, 3).Value = 43           # C24
$ws.Cells.Item(24, 4).Value = 59           # D24
$ws.Cells.Item(24, 5).Value = -27.118644067796   # E24
$ws.Cells.Item(24, 6).Value = 179          # F24
$ws.Cells.Item(24, 7).Value = 194          # G24
$ws.Cells.Item(24, 8).Value = -7.731958762886    # H24
$ws.Cells.Item(24, 9).Value = 1573         # I24
$ws.Cells.Item(24, 10).Value = 1733        # J24
$ws.Cells.Item(24, 11).Value = -9.232544720138   # K24
$ws.Cells.Item(24, 12).Value = 23.858267716535   # L24
$ws.Cells.Item(24, 13).Value = 19.257012888551   # M24

# ---------------------------------------------------------------------
# Row 25 — Misd. Assault
# ---------------------------------------------------------------------
$ws.Cells.Item(25, 3).Value = 11           # C25
$ws.Cells.Item(25, 4).Value = 11           # D25
$ws.Cells.Item(25, 5).Value = 0            # E25
$ws.Cells.Item(25, 6).Value = 59           # F25
$ws.Cells.Item(25, 7).Value = 45           # G25
$ws.Cells.Item(25, 8).Value = 31.111111111111    # H25
$ws.Cells.Item(25, 9).Value = 378          # I25
$ws.Cells.Item(25, 10).Value = 376         # J25
$ws.Cells.Item(25, 11).Value = 0.531914893617    # K25
$ws.Cells.Item(25, 12).Value = 24.752475247524   # L25
$ws.Cells.Item(25, 13).Value = 31.707317073170   # M25

# ---------------------------------------------------------------------
# Row 26 — UCR Rape*
# ---------------------------------------------------------------------
$ws.Cells.Item(26, 4).Value = 1            # D26
$ws.Cells.Item(26, 7).Value = 5            # G26
$ws.Cells.Item(26, 10).Value = 24          # J26
$ws.Cells.Item(26, 11).Value = -29.166666666666  # K26
$ws.Cells.Item(26, 12).Value = -29.166666666666  # L26

# ---------------------------------------------------------------------
# Row 27 — Other Sex Crimes
# ---------------------------------------------------------------------
$ws.Cells.Item(27, 3).Value = 2            # C27
$ws.Cells.Item(27, 4).Value = 5            # D27
$ws.Cells.Item(27, 5).Value = -60          # E27
$ws.Cells.Item(27, 6).Value = 8            # F27
$ws.Cells.Item(27, 7).Value = 10           # G27
$ws.Cells.Item(27, 8).Value = -20          # H27
$ws.Cells.Item(27, 9).Value = 81           # I27
$ws.Cells.Item(27, 10).Value = 74          # J27
$ws.Cells.Item(27, 11).Value = 9.459459459459   # K27
$ws.Cells.Item(27, 12).Value = -2.409638554216  # L27

# ---------------------------------------------------------------------
# Row 28 — Shooting Vic.  (C28, F28, I28 flip from numbers to "0" text)
# ---------------------------------------------------------------------
Set-TextFromNum 28 3 "0"        # C28 -> "0"
Set-TextFromNum 28 6 "0"        # F28 -> "0"
$ws.Cells.Item(28, 8).Value = -100         # H28
Set-TextFromNum 28 9 "0"        # I28 -> "0"
$ws.Cells.Item(28, 11).Value = -100        # K28
$ws.Cells.Item(28, 12).Value = -100        # L28
$ws.Cells.Item(28, 14).Value = -100        # N28

# ---------------------------------------------------------------------
# Row 29 — Shooting Inc.  (C29, F29, I29 flip from numbers to "0" text)
# ---------------------------------------------------------------------
Set-TextFromNum 29 3 "0"        # C29 -> "0"
Set-TextFromNum 29 6 "0"        # F29 -> "0"
$ws.Cells.Item(29, 8).Value = -100         # H29
Set-TextFromNum 29 9 "0"        # I29 -> "0"
$ws.Cells.Item(29, 11).Value = -100        # K29
$ws.Cells.Item(29, 12).Value = -100        # L29
$ws.Cells.Item(29, 14).Value = -100        # N29

# ---------------------------------------------------------------------
# Row 30 — Hate Crimes  (F30 flips from a number to "0" text)
# ---------------------------------------------------------------------
Set-TextFromNum 30 6 "0"        # F30 -> "0"
